# Rename the original sheet to "su" and build a second sheet "wu" as a
# sparser variant of the same report (same layout/labels, only row 6
# filled in with its own measurements).

$wb = $excel.ActiveWorkbook
$su = $wb.Worksheets.Item(1)
$su.Name = "su"

# Clone "su" (brings over labels, headers, merged cells, page setup, etc.)
# and place the clone right after it, then rename the clone to "wu".
$su.Copy($null, $su)
$wu = $wb.Worksheets.Item(2)
$wu.Name = "wu"

# --- Highlight row 6 on "su" with the new red/centered style and the
#     updated measurements (columns D,F,H,J,L,P + recalculated R6/S6). ---
$su.Range("A6:S6").Font.Color = 1972430
$su.Range("D6").Value = 1615
$su.Range("F6").Value = 1607
$su.Range("H6").Value = 1601
$su.Range("J6").Value = 1576
$su.Range("L6").Value = 1576
$su.Range("P6").Value = 1480

# --- "wu" only keeps the Name/#nodes/Original area columns (A:C) plus
#     the header rows; everything else outside row 6 is blanked out. ---
$wu.Range("D3:S5").ClearContents()
$wu.Range("D7:S14").ClearContents()
$wu.Range("B15:Q15").ClearContents()
$wu.Range("R15:S15").ClearContents()
$wu.Range("B15:Q15").HorizontalAlignment = -4108

# Row 6 on "wu" gets its own full set of measurements and the same
# red/centered highlight style as "su" row 6.
$wu.Range("A6:S6").Font.Color = 1972430
$wu.Range("D6").Value = 1615
$wu.Range("E6").Value = 52
$wu.Range("F6").Value = 1596
$wu.Range("G6").Value = 109
$wu.Range("H6").Value = 1582
$wu.Range("I6").Value = 139
$wu.Range("J6").Value = 1568
$wu.Range("K6").Value = 208
$wu.Range("L6").Value = 1562
$wu.Range("M6").Value = 230
$wu.Range("N6").Value = 1528
$wu.Range("O6").Value = 563
$wu.Range("P6").Value = 1487
$wu.Range("Q6").Value = 480

# --- View settings: "su" stays the active/selected tab, zoomed to 100%
#     with F6 selected; "wu" is left not-selected, zoomed to 100%, A6. ---
$su.Range("F6").Select()
$su.Activate()
